$wb = $excel.ActiveWorkbook

# The edited data lives on the "KW03" sheet, which is also the workbook's
# active tab.
$ws = $wb.Worksheets.Item("KW03")
$ws.Activate()

# Friday "time out" on the first week block changes from 18:00 to 12:00.
$ws.Range("F3").Value = 0.5

# Friday "time in" / "time out" on the second week block were empty and
# now get entered: 13:00 in, 20:00 out.
$ws.Range("F5").Value = 0.541666666666667
$ws.Range("F6").Value = 0.833333333333333

# Move the sheet's active-cell selection to F6.
$ws.Range("F6").Select()
